# Natmi following Dr Hou advice
# Update the NATMI LR-pair output (F8-Asgr2) numbers for rows 2-5 to reflect
# the re-run of the analysis (ligand/receptor expressing cell counts changed
# from 1 to 3, with corresponding recalculated expression/specificity values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.8546833333333335
$ws.Range("H2").Value = 2.56405
$ws.Range("I2").Value = 0.3097546281380014
$ws.Range("J2").Value = 0.3097546281380015
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.4418676666666667
$ws.Range("N2").Value = 1.325603
$ws.Range("Q2").Value = 0.377656930238889
$ws.Range("R2").Value = 3.398912372150001
$ws.Range("S2").Value = 0.3097546281380014
$ws.Range("T2").Value = 0.3097546281380015

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.251276
$ws.Range("H3").Value = 3.753828
$ws.Range("I3").Value = 0.453487879032787
$ws.Range("J3").Value = 0.4534878790327871
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.4418676666666667
$ws.Range("N3").Value = 1.325603
$ws.Range("Q3").Value = 0.5528984064760001
$ws.Range("R3").Value = 4.976085658284001
$ws.Range("S3").Value = 0.453487879032787
$ws.Range("T3").Value = 0.4534878790327871

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1689346666666667
$ws.Range("H4").Value = 0.506804
$ws.Range("I4").Value = 0.06122536009783416
$ws.Range("J4").Value = 0.06122536009783416
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4418676666666667
$ws.Range("N4").Value = 1.325603
$ws.Range("Q4").Value = 0.07464676697911113
$ws.Range("R4").Value = 0.6718209028120001
$ws.Range("S4").Value = 0.06122536009783416
$ws.Range("T4").Value = 0.06122536009783416

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.4843330000000001
$ws.Range("H5").Value = 1.452999
$ws.Range("I5").Value = 0.1755321327313773
$ws.Range("J5").Value = 0.1755321327313773
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4418676666666667
$ws.Range("N5").Value = 1.325603
$ws.Range("Q5").Value = 0.2140110925996667
$ws.Range("R5").Value = 1.926099833397
$ws.Range("S5").Value = 0.1755321327313773
$ws.Range("T5").Value = 0.1755321327313773
